$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 8334918
$ws.Range("J32").Value = 1763.2858
$ws.Range("L32").Value = 1763.2858
$ws.Range("N32").Value = -2415.2858

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 349.75
$ws.Range("I42").Value = 249.5
$ws.Range("K42").Value = 748.5
$ws.Range("M42").Value = -518.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3736
$ws.Range("I62").Value = 3729.6
$ws.Range("K62").Value = 3729.6
$ws.Range("M62").Value = -3105.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3736
$ws.Range("I65").Value = 3729.6
$ws.Range("K65").Value = 18648
$ws.Range("M65").Value = -15528

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3511.3157
$ws.Range("I137").Value = 1786.4546
$ws.Range("K137").Value = 5359.3638
$ws.Range("M137").Value = -2809.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1416.2898
$ws.Range("I32").Value = 1455
$ws.Range("K32").Value = 1455
$ws.Range("M32").Value = -1168

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 6900
$ws.Range("J21").Value = 6900
$ws.Range("L21").Value = 6900
$ws.Range("N21").Value = -7372

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 85000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 85000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 85000
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -86982

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3355.2856
$ws.Range("I99").Value = 2664.5
$ws.Range("K99").Value = 2664.5
$ws.Range("M99").Value = -1166.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 68565.39999999999
$ws.Range("J74").Value = 78206.75
$ws.Range("L74").Value = 78206.75
$ws.Range("N74").Value = -79954.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 68565.39999999999
$ws.Range("J77").Value = 78206.75
$ws.Range("L77").Value = 234620.25
$ws.Range("N77").Value = -243356.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 950.1539
$ws.Range("J94").Value = 1039.75
$ws.Range("L94").Value = 1039.75
$ws.Range("N94").Value = -1941.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3829
$ws.Range("I105").Value = 4075.75
$ws.Range("K105").Value = 4075.75
$ws.Range("M105").Value = -2328.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2708.1667
$ws.Range("I122").Value = 3066.75
$ws.Range("K122").Value = 9200.25
$ws.Range("M122").Value = -6750.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 223.04546
$ws.Range("I2").Value = 150.33333
$ws.Range("K2").Value = 901.9999799999999
$ws.Range("M2").Value = -788.9999799999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 9999.4
$ws.Range("I93").Value = 9999
$ws.Range("K93").Value = 29997
$ws.Range("M93").Value = -28125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 543.5
$ws.Range("I98").Value = 230.5
$ws.Range("J98").Value = 700
$ws.Range("K98").Value = 691.5
$ws.Range("L98").Value = 2100
$ws.Range("M98").Value = 806.5
$ws.Range("N98").Value = -5096

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 13904059
$ws.Range("I131").Value = 20854684
$ws.Range("K131").Value = 62564052
$ws.Range("M131").Value = -62559012

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6105.3184
$ws.Range("I126").Value = 1931.6154
$ws.Range("K126").Value = 5794.8462
$ws.Range("M126").Value = -3324.8462

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1578
$ws.Range("I7").Value = 1505.8667
$ws.Range("J7").Value = 1848.5
$ws.Range("K7").Value = 1505.8667
$ws.Range("L7").Value = 1848.5
$ws.Range("M7").Value = -1393.8667
$ws.Range("N7").Value = -2072.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9500
$ws.Range("J22").Value = 9000
$ws.Range("L22").Value = 9000
$ws.Range("N22").Value = -9590

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 9500
$ws.Range("J27").Value = 9000
$ws.Range("L27").Value = 9000
$ws.Range("N27").Value = -9214

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 26033
$ws.Range("I41").Value = 26033
$ws.Range("K41").Value = 26033
$ws.Range("M41").Value = -25595

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 35020.5
$ws.Range("I45").Value = 35020.5
$ws.Range("K45").Value = 35020.5
$ws.Range("M45").Value = -34613.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 32009.75
$ws.Range("I48").Value = 32009.75
$ws.Range("K48").Value = 32009.75
$ws.Range("M48").Value = -31348.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1098.6666
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1098.6666
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1098.6666
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -1820.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1098.6666
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1098.6666
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1098.6666
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -3594.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1432.7142
$ws.Range("I93").Value = 1135.75
$ws.Range("K93").Value = 1135.75
$ws.Range("M93").Value = 112.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 4627
$ws.Range("J100").Value = 5398.8
$ws.Range("L100").Value = 5398.8
$ws.Range("N100").Value = -6480.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1578
$ws.Range("I126").Value = 1505.8667
$ws.Range("J126").Value = 1848.5
$ws.Range("K126").Value = 4517.6001
$ws.Range("L126").Value = 5545.5
$ws.Range("M126").Value = -2047.6001
$ws.Range("N126").Value = -10485.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 12791.5
$ws.Range("I132").Value = 8999.75
$ws.Range("K132").Value = 26999.25
$ws.Range("M132").Value = -24469.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 55249.25
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 55249.25
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 55249.25
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -65449.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 57448
$ws.Range("I45").Value = 75000
$ws.Range("J45").Value = 53060
$ws.Range("K45").Value = 75000
$ws.Range("L45").Value = 53060
$ws.Range("M45").Value = -74509
$ws.Range("N45").Value = -54042

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3948.1667
$ws.Range("I96").Value = 1339.7142
$ws.Range("J96").Value = 7600
$ws.Range("K96").Value = 1339.7142
$ws.Range("L96").Value = 7600
$ws.Range("M96").Value = 33.28580000000011
$ws.Range("N96").Value = -10346

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2861.9
$ws.Range("I132").Value = 2577.5
$ws.Range("J132").Value = 3999.5
$ws.Range("K132").Value = 7732.5
$ws.Range("L132").Value = 11998.5
$ws.Range("M132").Value = -5202.5
$ws.Range("N132").Value = -17058.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1594521.9
$ws.Range("I136").Value = 8619.5
$ws.Range("K136").Value = 25858.5
$ws.Range("M136").Value = -23308.5

Write-Output "Applied all Gilgamesh_Profits updates"